$wb = $excel.ActiveWorkbook

# --- Update selection/active-cell state on the existing sheets ---

$wsZero = $wb.Worksheets.Item("ZERO_ROW_FOUR_COLUMN")
$wsZero.Range("D30").Select() | Out-Null

$wsNonBlank = $wb.Worksheets.Item("NON_BLANK_NO_HEADER")
$wsNonBlank.Range("M13").Select() | Out-Null

# --- Add the new worksheet at the end of the workbook ---

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "DISCON_3_COL_4_ROW"

# --- Populate the new sheet's data (disconnected 3-col/4-row block plus a
#     reference column E that documents the "discon" layout) ---

$newSheet.Range("A1").Value = "a"
$newSheet.Range("B1").Value = "b"
$newSheet.Range("C1").Value = "c"
$newSheet.Range("E1").Value = "discon"

$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = 2
$newSheet.Range("C2").Value = 3
$newSheet.Range("E2").Value = "a"

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = 1
$newSheet.Range("C3").Value = 1
$newSheet.Range("E3").Value = 1

$newSheet.Range("A4").Value = 1
$newSheet.Range("B4").Value = 1
$newSheet.Range("C4").Value = 1
$newSheet.Range("E4").Value = 1

$newSheet.Range("A5").Value = 1
$newSheet.Range("B5").Value = 1
$newSheet.Range("C5").Value = 1
$newSheet.Range("D5").Value = 1
$newSheet.Range("E5").Value = 1

$newSheet.Range("D6").Value = 1
$newSheet.Range("E6").Value = 1

# Final selection on the new (now active) sheet
$newSheet.Range("A7").Select() | Out-Null
